# Generate Report for Archive
#
# On the "Overview" sheet, rows for files 59847de3-....md and
# 66e28f57-....md move from "Ready for handoff" to "In Translation"
# for both the zh-cn (column E) and de-de (column F) languages.
# The f24b7175-....md row (and everything else) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"

$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
